# Apply commit "Add data for 2022-09-02":
# The report's "as of" cutoff date moves from August 24 to August 25,
# so the current-year column header/sheet name shift by one day, and a
# handful of neighborhood/month cells (the "August" column for every
# year, since the cutoff day-of-year affects all years' August totals)
# get bumped up by the newly-included day's carjacking counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the "August 2022 (through August NN)" label.
$ws.Name = "Through 2022-08-25"
$ws.Range("B1").Value = "August 2022 (through August 25)"

# Cells whose totals increase by the events that occurred on Aug 25 of
# the respective year (existing values bumped by 1).
$ws.Range("R2").Value = 8
$ws.Range("AX2").Value = 3
$ws.Range("AP3").Value = 4
$ws.Range("B4").Value = 5
$ws.Range("R4").Value = 3
$ws.Range("J5").Value = 13
$ws.Range("AX5").Value = 5
$ws.Range("R7").Value = 6
$ws.Range("AX12").Value = 3
$ws.Range("Z13").Value = 4
$ws.Range("B27").Value = 2
$ws.Range("AX66").Value = 4
$ws.Range("AH83").Value = 2
$ws.Range("AH92").Value = 3

# Cells that previously had no recorded carjackings for that
# neighborhood/month and now have exactly 1 (newly populated cells).
$ws.Range("AP10").Value = 1
$ws.Range("AP11").Value = 1
$ws.Range("Z23").Value = 1
$ws.Range("J25").Value = 1
$ws.Range("J43").Value = 1
$ws.Range("B44").Value = 1
